$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/issue number and report date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "49"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "12/5/2022"
$c9.Characters(47, 9).Text = "12/11/2022"

# --- Weekly crime statistics table (rows 14-30, columns C-N) ---
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 9
$ws.Range("E14").Value = 11.111111111111
$ws.Range("F14").Value = 29
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = -3.333333333333
$ws.Range("I14").Value = 405
$ws.Range("J14").Value = 457
$ws.Range("K14").Value = -11.378555798687
$ws.Range("L14").Value = -8.371040723981
$ws.Range("M14").Value = -21.052631578947
$ws.Range("N14").Value = -77.747252747252
$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 29
$ws.Range("E15").Value = -31.034482758620
$ws.Range("F15").Value = 92
$ws.Range("G15").Value = 119
$ws.Range("H15").Value = -22.689075630252
$ws.Range("I15").Value = 1535
$ws.Range("J15").Value = 1424
$ws.Range("K15").Value = 7.794943820224
$ws.Range("L15").Value = 10.830324909747
$ws.Range("M15").Value = 16.464339908953
$ws.Range("N15").Value = -50.403877221324
$ws.Range("C16").Value = 285
$ws.Range("D16").Value = 307
$ws.Range("E16").Value = -7.166123778501
$ws.Range("F16").Value = 1251
$ws.Range("G16").Value = 1310
$ws.Range("H16").Value = -4.503816793893
$ws.Range("I16").Value = 16576
$ws.Range("J16").Value = 13000
$ws.Range("K16").Value = 27.507692307692
$ws.Range("L16").Value = 33.958299660578
$ws.Range("M16").Value = -10.254466702761
$ws.Range("N16").Value = -79.606048302759
$ws.Range("C17").Value = 444
$ws.Range("D17").Value = 465
$ws.Range("E17").Value = -4.516129032258
$ws.Range("F17").Value = 1780
$ws.Range("G17").Value = 1823
$ws.Range("H17").Value = -2.358749314317
$ws.Range("I17").Value = 24690
$ws.Range("J17").Value = 21896
$ws.Range("K17").Value = 12.760321519912
$ws.Range("L17").Value = 25.031650377272
$ws.Range("M17").Value = 51.677110210099
$ws.Range("N17").Value = -37.277715679300
$ws.Range("C18").Value = 292
$ws.Range("D18").Value = 306
$ws.Range("E18").Value = -4.575163398692
$ws.Range("F18").Value = 1166
$ws.Range("G18").Value = 1253
$ws.Range("H18").Value = -6.943335993615
$ws.Range("I18").Value = 14909
$ws.Range("J18").Value = 11981
$ws.Range("K18").Value = 24.438694599783
$ws.Range("L18").Value = 1.256452051073
$ws.Range("M18").Value = -16.048200912213
$ws.Range("N18").Value = -84.397304142160
$ws.Range("C19").Value = 987
$ws.Range("D19").Value = 1416
$ws.Range("E19").Value = -30.296610169491
$ws.Range("F19").Value = 3830
$ws.Range("G19").Value = 4821
$ws.Range("H19").Value = -20.555901265297
$ws.Range("I19").Value = 48859
$ws.Range("J19").Value = 37746
$ws.Range("K19").Value = 29.441530228368
$ws.Range("L19").Value = 43.774828590766
$ws.Range("M19").Value = 35.734526058450
$ws.Range("N19").Value = -39.876945794622
$ws.Range("C20").Value = 301
$ws.Range("D20").Value = 208
$ws.Range("E20").Value = 44.711538461538
$ws.Range("F20").Value = 1105
$ws.Range("G20").Value = 939
$ws.Range("H20").Value = 17.678381256656
$ws.Range("I20").Value = 12901
$ws.Range("J20").Value = 9757
$ws.Range("K20").Value = 32.223019370708
$ws.Range("L20").Value = 49.959316517493
$ws.Range("M20").Value = 31.535481239804
$ws.Range("N20").Value = -87.849763135836
$ws.Range("C21").Value = 2339
$ws.Range("D21").Value = 2740
$ws.Range("E21").Value = -14.635036496350
$ws.Range("F21").Value = 9253
$ws.Range("G21").Value = 10295
$ws.Range("H21").Value = -10.121418164157
$ws.Range("I21").Value = 119875
$ws.Range("J21").Value = 96261
$ws.Range("K21").Value = 24.531222405751
$ws.Range("L21").Value = 31.358346665497
$ws.Range("M21").Value = 19.705018873200
$ws.Range("N21").Value = -70.658857047748
$ws.Range("C22").Value = 44
$ws.Range("D22").Value = 55
$ws.Range("E22").Value = -20
$ws.Range("F22").Value = 169
$ws.Range("G22").Value = 221
$ws.Range("H22").Value = -23.529411764705
$ws.Range("I22").Value = 2192
$ws.Range("J22").Value = 1684
$ws.Range("K22").Value = 30.166270783848
$ws.Range("L22").Value = 30.321046373365
$ws.Range("M22").Value = 6.926829268292
$ws.Range("C23").Value = 96
$ws.Range("D23").Value = 118
$ws.Range("E23").Value = -18.644067796610
$ws.Range("F23").Value = 441
$ws.Range("G23").Value = 480
$ws.Range("H23").Value = -8.125
$ws.Range("I23").Value = 5634
$ws.Range("J23").Value = 5255
$ws.Range("K23").Value = 7.212178877259
$ws.Range("L23").Value = 15.261865793780
$ws.Range("M23").Value = 40.603943099575
$ws.Range("C24").Value = 2317
$ws.Range("D24").Value = 2055
$ws.Range("E24").Value = 12.749391727493
$ws.Range("F24").Value = 9059
$ws.Range("G24").Value = 8149
$ws.Range("H24").Value = 11.167014357589
$ws.Range("I24").Value = 109563
$ws.Range("J24").Value = 81419
$ws.Range("K24").Value = 34.566870140876
$ws.Range("L24").Value = 41.417231364956
$ws.Range("M24").Value = 41.280464216634
$ws.Range("C25").Value = 733
$ws.Range("D25").Value = 790
$ws.Range("E25").Value = -7.215189873417
$ws.Range("F25").Value = 2880
$ws.Range("G25").Value = 3021
$ws.Range("H25").Value = -4.667328699106
$ws.Range("I25").Value = 39169
$ws.Range("J25").Value = 34407
$ws.Range("K25").Value = 13.840206934635
$ws.Range("L25").Value = 23.932922006011
$ws.Range("M25").Value = -10.125739984397
$ws.Range("C26").Value = 39
$ws.Range("D26").Value = 39
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 159
$ws.Range("G26").Value = 197
$ws.Range("H26").Value = -19.289340101522
$ws.Range("I26").Value = 2468
$ws.Range("J26").Value = 2340
$ws.Range("K26").Value = 5.470085470085
$ws.Range("L26").Value = 14.259259259259
$ws.Range("C27").Value = 98
$ws.Range("D27").Value = 121
$ws.Range("E27").Value = -19.008264462809
$ws.Range("F27").Value = 377
$ws.Range("G27").Value = 416
$ws.Range("H27").Value = -9.375
$ws.Range("I27").Value = 4928
$ws.Range("J27").Value = 4699
$ws.Range("K27").Value = 4.873377314322
$ws.Range("L27").Value = 35.720187276232
$ws.Range("C28").Value = 28
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = -6.666666666666
$ws.Range("F28").Value = 80
$ws.Range("G28").Value = 136
$ws.Range("H28").Value = -41.176470588235
$ws.Range("I28").Value = 1502
$ws.Range("J28").Value = 1787
$ws.Range("K28").Value = -15.948517067711
$ws.Range("L28").Value = -15.759955131800
$ws.Range("M28").Value = -11.281748375664
$ws.Range("N28").Value = -73.221608129791
$ws.Range("C29").Value = 22
$ws.Range("D29").Value = 23
$ws.Range("E29").Value = -4.347826086956
$ws.Range("F29").Value = 71
$ws.Range("G29").Value = 112
$ws.Range("H29").Value = -36.607142857142
$ws.Range("I29").Value = 1240
$ws.Range("J29").Value = 1490
$ws.Range("K29").Value = -16.778523489932
$ws.Range("L29").Value = -14.951989026063
$ws.Range("M29").Value = -11.428571428571
$ws.Range("N29").Value = -75.391942845802
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = -70
$ws.Range("F30").Value = 20
$ws.Range("G30").Value = 30
$ws.Range("H30").Value = -33.333333333333
$ws.Range("I30").Value = 578
$ws.Range("J30").Value = 506
$ws.Range("K30").Value = 14.229249011857
$ws.Range("L30").Value = 131.2

# --- Column G auto-fit to reflect new widest value (bestFit) ---
$ws.Columns.Item(7).AutoFit()
